$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Shift the data in rows 2-13 one column to the right (A:E -> B:F),
# processing from the rightmost column first so we don't clobber data
# we still need to read. The header row (row 1) is left untouched.
for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 5; $c -ge 1; $c--) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($r, $c + 1)
        $dstCell.ClearContents()
        $srcCell.Copy($dstCell)
    }
}

# Add header for the new "Link" column, matching the style of the other headers.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 6))
$ws.Cells.Item(1, 6).Value = "Link"

# Fill column A (rows 2-13) with the CVE id extracted from the link URL (now in column F)
for ($r = 2; $r -le $lastRow; $r++) {
    $link = $ws.Cells.Item($r, 6).Value()
    $cve = $link.Substring($link.LastIndexOf("/") + 1)
    $ws.Cells.Item($r, 1).Value = $cve
}
